$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text first (shared string change DNI -> DOCUMENTO)
$ws.Range("B2").Value = "DOCUMENTO "

# Move the header cell (value + formatting) from B2 to A1
$ws.Range("B2").Cut($ws.Range("A1")) | Out-Null

# B2 is now vacated; remove any leftover formatting/content so it goes back to default
$ws.Range("B2").Clear()

$wb.Save()
